# Updated cryptos list on Mon Sep  4 17:51:44 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.924.02'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '1.633.95'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.08'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5094'
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2579'
$ws.Range("E8").Value = '  +1.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06351'
$ws.Range("E9").Value = '  +0.40%  '
$ws.Range("E10").Value = '  +0.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07790'
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.270'
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").Value = '1.639.46'
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("D14").Value = '1.859.43'
$ws.Range("E14").Value = '  +0.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5511'
$ws.Range("E15").Value = '  +1.93%  '
$ws.Range("D17").Value = '0.0₅7668'
$ws.Range("E17").Value = '  -0.39%  '
$ws.Range("D18").Value = '25.939.57'
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '196.04'
$ws.Range("E20").Value = '  +0.54%  '
$ws.Range("E21").Value = '  +0.27%  '
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.065'
$ws.Range("E23").Value = '  +1.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.910'
$ws.Range("E25").Value = '  +2.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.24'
$ws.Range("E26").Value = '  +1.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1249'
$ws.Range("E27").Value = '  +4.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.63'
$ws.Range("E28").Value = '  +0.44%  '
$ws.Range("E29").Value = '  -0.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.241'
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04911'
$ws.Range("E31").Value = '  +0.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.248'
$ws.Range("E32").Value = '  +0.43%  '
$ws.Range("E33").Value = '  +1.28%  '
$ws.Range("E34").Value = '  +1.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.370'
$ws.Range("E35").Value = '  +0.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8986'
$ws.Range("E36").Value = '  +1.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5539'
$ws.Range("E37").Value = '  +2.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.536'
$ws.Range("E38").Value = '  -1.38%  '
$ws.Range("D39").Value = '1.113.38'
$ws.Range("E39").Value = '  -2.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01558'
$ws.Range("E40").Value = '  +0.86%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.619'
$ws.Range("E42").Value = '  +3.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7962'
$ws.Range("E43").Value = '  -1.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.50'
$ws.Range("E45").Value = '  -5.13%  '
$ws.Range("D46").Value = '1.769.49'
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4448'
$ws.Range("E47").Value = '  -1.63%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.004'
$ws.Range("E48").Value = '  +0.56%  '
$ws.Range("E49").Value = '  +0.45%  '
$ws.Range("E50").Value = '  +1.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.575'
$ws.Range("E51").Value = '  +4.01%  '
